# Apply odds updates to rows 3 and 4 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    # Row 3
    "G3"  = 1.22
    "H3"  = 6
    "I3"  = 12
    "J3"  = 1.62
    "K3"  = 2.88
    "L3"  = 8.5
    "Q3"  = 1.44
    "R3"  = 2.7
    "S3"  = 1.22
    "T3"  = 4
    "U3"  = 1.83
    "V3"  = 1.83
    "Z3"  = 8
    "AA3" = 10
    "AD3" = 12
    "AE3" = 21
    "AG3" = 34
    "AI3" = 34
    "AJ3" = 151
    "AK3" = 67
    "AM3" = 251
    "AN3" = 3.4
    "AO3" = 5.5
    "AQ3" = 13
    "AT3" = 4
    "AU3" = 9.5
    "AW3" = 11
    "AY3" = 41
    "AZ3" = 201
    "BA3" = 151
    "BB3" = 251

    # Row 4
    "L4"  = 6
    "O4"  = 1.5
    "P4"  = 2.63
    "Q4"  = 2.5
    "R4"  = 1.5
    "Z4"  = 13
    "AI4" = 19
    "BB4" = 501
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
